$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (column C) date value for rows 2-6 from 45224 to 45233
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
